$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Log_Muestras")

$timestamps = @{
    2 = "2025-10-17T07:09:27.427522"
    3 = "2025-10-17T07:09:27.427522"
    4 = "2025-10-17T07:09:27.427522"
    5 = "2025-10-17T07:09:27.430114"
    6 = "2025-10-17T07:09:27.430114"
    7 = "2025-10-17T07:09:27.430637"
    8 = "2025-10-17T07:09:27.430637"
    9 = "2025-10-17T07:09:27.430637"
    10 = "2025-10-17T07:09:27.430637"
    11 = "2025-10-17T07:09:27.431737"
    12 = "2025-10-17T07:09:27.431737"
    13 = "2025-10-17T07:09:27.432250"
    14 = "2025-10-17T07:09:27.432250"
    15 = "2025-10-17T07:09:27.432250"
    16 = "2025-10-17T07:09:27.491557"
    17 = "2025-10-17T07:09:27.492557"
    18 = "2025-10-17T07:09:27.492557"
    19 = "2025-10-17T07:09:27.492557"
    20 = "2025-10-17T07:09:27.492557"
    21 = "2025-10-17T07:09:27.494129"
    22 = "2025-10-17T07:09:27.494530"
    23 = "2025-10-17T07:09:27.494530"
    24 = "2025-10-17T07:09:27.494530"
    25 = "2025-10-17T07:09:27.494530"
    26 = "2025-10-17T07:09:27.564640"
    27 = "2025-10-17T07:09:27.564640"
    28 = "2025-10-17T07:09:27.564640"
    29 = "2025-10-17T07:09:27.565640"
    30 = "2025-10-17T07:09:27.565640"
    31 = "2025-10-17T07:09:27.565640"
    32 = "2025-10-17T07:09:27.565640"
    33 = "2025-10-17T07:09:27.565640"
    34 = "2025-10-17T07:09:27.565640"
    35 = "2025-10-17T07:09:27.565640"
    36 = "2025-10-17T07:09:27.565640"
    37 = "2025-10-17T07:09:27.565640"
    38 = "2025-10-17T07:09:27.565640"
    39 = "2025-10-17T07:09:27.565640"
    40 = "2025-10-17T07:09:27.566640"
    41 = "2025-10-17T07:09:27.566640"
    42 = "2025-10-17T07:09:27.566640"
    43 = "2025-10-17T07:09:27.566640"
    44 = "2025-10-17T07:09:27.566640"
    45 = "2025-10-17T07:09:27.566640"
    46 = "2025-10-17T07:09:27.566640"
    47 = "2025-10-17T07:09:27.566640"
    48 = "2025-10-17T07:09:27.566640"
}

foreach ($row in $timestamps.Keys) {
    $ws.Cells.Item($row, 26).Value = $timestamps[$row]
}

Write-Output "Updated $($timestamps.Count) timestamp cells in column Z"
